$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F column (想去人数 / "want to go" counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6850
$ws1.Range("F4").Value = 439
$ws1.Range("F5").Value = 73
$ws1.Range("F8").Value = 112
$ws1.Range("F12").Value = 37
$ws1.Range("F13").Value = 183
$ws1.Range("F14").Value = 429
$ws1.Range("F16").Value = 1785
$ws1.Range("F17").Value = 29
$ws1.Range("F18").Value = 3467
$ws1.Range("F20").Value = 237
$ws1.Range("F21").Value = 16
$ws1.Range("F22").Value = 2102
$ws1.Range("F23").Value = 189
$ws1.Range("F29").Value = 140

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5

# Sheet "全部类型" (All types) - mirrors "展览" but with an extra row (row 7) from "演出"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6850
$ws4.Range("F4").Value = 439
$ws4.Range("F5").Value = 73
$ws4.Range("F7").Value = 5
$ws4.Range("F9").Value = 112
$ws4.Range("F13").Value = 37
$ws4.Range("F14").Value = 183
$ws4.Range("F15").Value = 429
$ws4.Range("F17").Value = 1785
$ws4.Range("F18").Value = 29
$ws4.Range("F19").Value = 3467
$ws4.Range("F21").Value = 237
$ws4.Range("F22").Value = 16
$ws4.Range("F23").Value = 2102
$ws4.Range("F24").Value = 189
$ws4.Range("F30").Value = 140
